$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 633.0714
$ws.Range("J2").Value = 877.875
$ws.Range("L2").Value = 877.875
$ws.Range("N2").Value = -1103.875
$ws.Range("H12").Value = 4884.4546
$ws.Range("I12").Value = 5608.737
$ws.Range("J12").Value = 297.33334
$ws.Range("K12").Value = 5608.737
$ws.Range("L12").Value = 297.33334
$ws.Range("M12").Value = -5438.737
$ws.Range("N12").Value = -637.33334
$ws.Range("H17").Value = 414173.94
$ws.Range("J17").Value = 414173.94
$ws.Range("L17").Value = 1242521.82
$ws.Range("N17").Value = -1242857.82
$ws.Range("H32").Value = 7888.0527
$ws.Range("I32").Value = 7677.9
$ws.Range("J32").Value = 8121.5557
$ws.Range("K32").Value = 7677.9
$ws.Range("L32").Value = 8121.5557
$ws.Range("M32").Value = -7351.9
$ws.Range("N32").Value = -8773.555700000001
$ws.Range("H53").Value = 278
$ws.Range("J53").Value = 331
$ws.Range("L53").Value = 331
$ws.Range("N53").Value = -1605
$ws.Range("H62").Value = 4545.727
$ws.Range("J62").Value = 5599.6
$ws.Range("L62").Value = 5599.6
$ws.Range("N62").Value = -6847.6
$ws.Range("H65").Value = 4545.727
$ws.Range("J65").Value = 5599.6
$ws.Range("L65").Value = 27998
$ws.Range("N65").Value = -34238
$ws.Range("H70").Value = 10504.706
$ws.Range("J70").Value = 11067.5
$ws.Range("L70").Value = 33202.5
$ws.Range("N70").Value = -33742.5
$ws.Range("H73").Value = 10504.706
$ws.Range("J73").Value = 11067.5
$ws.Range("L73").Value = 33202.5
$ws.Range("N73").Value = -35074.5
$ws.Range("H75").Value = 43500
$ws.Range("J75").Value = 43500
$ws.Range("L75").Value = 43500
$ws.Range("N75").Value = -45372
$ws.Range("H76").Value = 5625.75
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684
$ws.Range("H78").Value = 43500
$ws.Range("J78").Value = 43500
$ws.Range("L78").Value = 130500
$ws.Range("N78").Value = -139860
$ws.Range("H79").Value = 5625.75
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907
$ws.Range("H80").Value = 344.55554
$ws.Range("I80").Value = 251.18182
$ws.Range("J80").Value = 491.2857
$ws.Range("K80").Value = 753.5454599999999
$ws.Range("L80").Value = 1473.8571
$ws.Range("M80").Value = 244.4545400000001
$ws.Range("N80").Value = -3469.8571
$ws.Range("H83").Value = 344.55554
$ws.Range("I83").Value = 251.18182
$ws.Range("J83").Value = 491.2857
$ws.Range("K83").Value = 2260.63638
$ws.Range("L83").Value = 4421.571300000001
$ws.Range("M83").Value = 2731.36362
$ws.Range("N83").Value = -14405.5713
$ws.Range("H86").Value = 10432.214
$ws.Range("I86").Value = 11700.375
$ws.Range("K86").Value = 11700.375
$ws.Range("M86").Value = -10577.375
$ws.Range("H89").Value = 10432.214
$ws.Range("I89").Value = 11700.375
$ws.Range("K89").Value = 58501.875
$ws.Range("M89").Value = -52885.875
$ws.Range("H103").Value = 991.8
$ws.Range("J103").Value = 991.8
$ws.Range("L103").Value = 2975.4
$ws.Range("N103").Value = -4147.4
$ws.Range("H116").Value = 5258.654
$ws.Range("I116").Value = 3457.2222
$ws.Range("J116").Value = 9311.875
$ws.Range("K116").Value = 3457.2222
$ws.Range("L116").Value = 9311.875
$ws.Range("M116").Value = -15.22220000000016
$ws.Range("N116").Value = -16195.875
$ws.Range("H121").Value = 1849.75
$ws.Range("J121").Value = 1849.75
$ws.Range("L121").Value = 5549.25
$ws.Range("N121").Value = -9043.25
$ws.Range("H132").Value = 2397.0293
$ws.Range("I132").Value = 2397.0293
$ws.Range("K132").Value = 7191.0879
$ws.Range("M132").Value = -4661.0879
$ws.Range("H137").Value = 1824.3572
$ws.Range("I137").Value = 1321.909
$ws.Range("J137").Value = 3666.6667
$ws.Range("K137").Value = 3965.727
$ws.Range("L137").Value = 11000.0001
$ws.Range("M137").Value = -1415.727
$ws.Range("N137").Value = -16100.0001
# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 924904.4399999999
$ws.Range("I2").Value = 1055962.2
$ws.Range("K2").Value = 1055962.2
$ws.Range("M2").Value = -1055849.2
$ws.Range("H5").Value = 905.2857
$ws.Range("I5").Value = 821
$ws.Range("K5").Value = 821
$ws.Range("M5").Value = -709
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""
$ws.Range("H28").Value = 21499.75
$ws.Range("I28").Value = 15333.333
$ws.Range("J28").Value = 39999
$ws.Range("K28").Value = 15333.333
$ws.Range("L28").Value = 39999
$ws.Range("M28").Value = -15141.333
$ws.Range("N28").Value = -40383
$ws.Range("H32").Value = 1192.9711
$ws.Range("I32").Value = 1179.0151
$ws.Range("K32").Value = 1179.0151
$ws.Range("M32").Value = -892.0151000000001
$ws.Range("H40").Value = 20028
$ws.Range("I40").Value = 20028
$ws.Range("K40").Value = 20028
$ws.Range("M40").Value = -19852
$ws.Range("H45").Value = 2177.6
$ws.Range("I45").Value = 1972
$ws.Range("K45").Value = 1972
$ws.Range("M45").Value = -1595
$ws.Range("H63").Value = 2343.9333
$ws.Range("J63").Value = 1239.5
$ws.Range("L63").Value = 1239.5
$ws.Range("N63").Value = -2611.5
$ws.Range("H66").Value = 2343.9333
$ws.Range("J66").Value = 1239.5
$ws.Range("L66").Value = 6197.5
$ws.Range("N66").Value = -13061.5
$ws.Range("H74").Value = 43484732
$ws.Range("I74").Value = 50006092
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 50006092
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -50005218
$ws.Range("N74").Value = -10748
$ws.Range("H76").Value = 34000
$ws.Range("J76").Value = 34000
$ws.Range("L76").Value = 34000
$ws.Range("N76").Value = -34676
$ws.Range("H77").Value = 43484732
$ws.Range("I77").Value = 50006092
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 250030460
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -250026092
$ws.Range("N77").Value = -53736
$ws.Range("H79").Value = 34000
$ws.Range("J79").Value = 34000
$ws.Range("L79").Value = 34000
$ws.Range("N79").Value = -36340
$ws.Range("H99").Value = 21499.75
$ws.Range("I99").Value = 15333.333
$ws.Range("J99").Value = 39999
$ws.Range("K99").Value = 15333.333
$ws.Range("L99").Value = 39999
$ws.Range("M99").Value = -12338.333
$ws.Range("N99").Value = -45989
$ws.Range("H110").Value = 146281.14
$ws.Range("I110").Value = 169661.5
$ws.Range("J110").Value = 5999
$ws.Range("K110").Value = 169661.5
$ws.Range("L110").Value = 5999
$ws.Range("M110").Value = -167616.5
$ws.Range("N110").Value = -10089
$ws.Range("H116").Value = 924904.4399999999
$ws.Range("I116").Value = 1055962.2
$ws.Range("K116").Value = 1055962.2
$ws.Range("M116").Value = -1053668.2
$ws.Range("H132").Value = 8338723.5
$ws.Range("I132").Value = 11115434
$ws.Range("J132").Value = 8591.666999999999
$ws.Range("K132").Value = 33346302
$ws.Range("L132").Value = 25775.001
$ws.Range("M132").Value = -33343772
$ws.Range("N132").Value = -30835.001
# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 924904.4399999999
$ws.Range("I3").Value = 1055962.2
$ws.Range("K3").Value = 1055962.2
$ws.Range("M3").Value = -1055848.2
$ws.Range("H4").Value = 905.2857
$ws.Range("I4").Value = 821
$ws.Range("K4").Value = 821
$ws.Range("M4").Value = -706
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("H86").Value = 4000.5
$ws.Range("I86").Value = 4000.5
$ws.Range("K86").Value = 4000.5
$ws.Range("M86").Value = -2877.5
$ws.Range("H89").Value = 4000.5
$ws.Range("I89").Value = 4000.5
$ws.Range("K89").Value = 20002.5
$ws.Range("M89").Value = -14386.5
$ws.Range("H99").Value = 1992.7878
$ws.Range("I99").Value = 1810.2667
$ws.Range("K99").Value = 1810.2667
$ws.Range("M99").Value = -312.2666999999999
$ws.Range("H107").Value = 128496.25
$ws.Range("I107").Value = 3598.2
$ws.Range("K107").Value = 3598.2
$ws.Range("M107").Value = -1678.2
$ws.Range("H138").Value = 125000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 125000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 125000
$ws.Range("M138").Value = ""
$ws.Range("N138").Value = -135280
# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 207.14285
$ws.Range("I7").Value = 60
$ws.Range("J7").Value = 266
$ws.Range("K7").Value = 60
$ws.Range("L7").Value = 266
$ws.Range("M7").Value = 53
$ws.Range("N7").Value = -492
$ws.Range("H16").Value = 6699.75
$ws.Range("I16").Value = 1800
$ws.Range("K16").Value = 1800
$ws.Range("M16").Value = -1513
$ws.Range("H31").Value = 6072.0586
$ws.Range("J31").Value = 6997.364
$ws.Range("L31").Value = 6997.364
$ws.Range("N31").Value = -7587.364
$ws.Range("H34").Value = 6072.0586
$ws.Range("J34").Value = 6997.364
$ws.Range("L34").Value = 6997.364
$ws.Range("N34").Value = -7401.364
$ws.Range("H52").Value = 58999.5
$ws.Range("J52").Value = 60000
$ws.Range("L52").Value = 60000
$ws.Range("N52").Value = -60588
$ws.Range("H58").Value = 20838000
$ws.Range("I58").Value = 21743988
$ws.Range("J58").Value = 295
$ws.Range("K58").Value = 21743988
$ws.Range("L58").Value = 295
$ws.Range("M58").Value = -21743785
$ws.Range("N58").Value = -701
$ws.Range("H62").Value = 4197.875
$ws.Range("J62").Value = 4600
$ws.Range("L62").Value = 4600
$ws.Range("N62").Value = -5848
$ws.Range("H65").Value = 4197.875
$ws.Range("J65").Value = 4600
$ws.Range("L65").Value = 23000
$ws.Range("N65").Value = -29240
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101498
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -307488
$ws.Range("H82").Value = 50000
$ws.Range("I82").Value = 50000
$ws.Range("K82").Value = 50000
$ws.Range("M82").Value = -49639
$ws.Range("H85").Value = 50000
$ws.Range("I85").Value = 50000
$ws.Range("K85").Value = 50000
$ws.Range("M85").Value = -48752
$ws.Range("H86").Value = 5004.4287
$ws.Range("I86").Value = 5164.3335
$ws.Range("J86").Value = 4716.6
$ws.Range("K86").Value = 5164.3335
$ws.Range("L86").Value = 4716.6
$ws.Range("M86").Value = -4041.3335
$ws.Range("N86").Value = -6962.6
$ws.Range("H89").Value = 5004.4287
$ws.Range("I89").Value = 5164.3335
$ws.Range("J89").Value = 4716.6
$ws.Range("K89").Value = 25821.6675
$ws.Range("L89").Value = 23583
$ws.Range("M89").Value = -20205.6675
$ws.Range("N89").Value = -34815
$ws.Range("H99").Value = 4635
$ws.Range("I99").Value = 4635
$ws.Range("K99").Value = 4635
$ws.Range("M99").Value = -3137
$ws.Range("H113").Value = 6699.75
$ws.Range("I113").Value = 1800
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370
$ws.Range("H126").Value = 4635
$ws.Range("I126").Value = 4635
$ws.Range("K126").Value = 13905
$ws.Range("M126").Value = -11435
$ws.Range("H132").Value = 21741004
$ws.Range("I132").Value = 27779294
$ws.Range("J132").Value = 3162.4
$ws.Range("K132").Value = 83337882
$ws.Range("L132").Value = 9487.200000000001
$ws.Range("M132").Value = -83335352
$ws.Range("N132").Value = -14547.2
$ws.Range("H134").Value = 6580280.5
$ws.Range("I134").Value = 6945792
$ws.Range("K134").Value = 20837376
$ws.Range("M134").Value = -20834841
$ws.Range("H136").Value = 20838000
$ws.Range("I136").Value = 21743988
$ws.Range("J136").Value = 295
$ws.Range("K136").Value = 65231964
$ws.Range("L136").Value = 885
$ws.Range("M136").Value = -65229414
$ws.Range("N136").Value = -5985
# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H11").Value = 157213.84
$ws.Range("I11").Value = 165892.28
$ws.Range("J11").Value = 1002
$ws.Range("K11").Value = 497676.84
$ws.Range("L11").Value = 3006
$ws.Range("M11").Value = -497536.84
$ws.Range("N11").Value = -3286
$ws.Range("H14").Value = 357.44446
$ws.Range("I14").Value = 357.44446
$ws.Range("K14").Value = 1072.33338
$ws.Range("M14").Value = -899.33338
$ws.Range("H34").Value = 912.7143
$ws.Range("J34").Value = 945
$ws.Range("L34").Value = 2835
$ws.Range("N34").Value = -3003
$ws.Range("H39").Value = 2831.6
$ws.Range("I39").Value = 993.2
$ws.Range("J39").Value = 4670
$ws.Range("K39").Value = 2979.6
$ws.Range("L39").Value = 14010
$ws.Range("M39").Value = -2685.6
$ws.Range("N39").Value = -14598
$ws.Range("H55").Value = 2520.4
$ws.Range("I55").Value = 704
$ws.Range("J55").Value = 2722.2222
$ws.Range("K55").Value = 2112
$ws.Range("L55").Value = 8166.6666
$ws.Range("M55").Value = -1935
$ws.Range("N55").Value = -8520.6666
$ws.Range("H56").Value = 242380.56
$ws.Range("I56").Value = 242380.56
$ws.Range("K56").Value = 242380.56
$ws.Range("M56").Value = -241850.56
$ws.Range("H64").Value = 917343.25
$ws.Range("J64").Value = 10080
$ws.Range("L64").Value = 30240
$ws.Range("N64").Value = -30780
$ws.Range("H67").Value = 917343.25
$ws.Range("J67").Value = 10080
$ws.Range("L67").Value = 30240
$ws.Range("N67").Value = -32112
$ws.Range("H68").Value = 1599.3334
$ws.Range("J68").Value = 1599.3334
$ws.Range("L68").Value = 4798.0002
$ws.Range("N68").Value = -6420.0002
$ws.Range("H70").Value = 16925
$ws.Range("I70").Value = 7700
$ws.Range("K70").Value = 23100
$ws.Range("M70").Value = -22785
$ws.Range("H71").Value = 1599.3334
$ws.Range("J71").Value = 1599.3334
$ws.Range("L71").Value = 14394.0006
$ws.Range("N71").Value = -22506.0006
$ws.Range("H73").Value = 16925
$ws.Range("I73").Value = 7700
$ws.Range("K73").Value = 23100
$ws.Range("M73").Value = -22008
$ws.Range("H92").Value = 399.5
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H97").Value = 1060.4117
$ws.Range("I97").Value = 609.6
$ws.Range("K97").Value = 1828.8
$ws.Range("M97").Value = -1332.8
$ws.Range("H107").Value = 890.3125
$ws.Range("J107").Value = 1069.8
$ws.Range("L107").Value = 3209.4
$ws.Range("N107").Value = -7049.4
$ws.Range("H115").Value = 9999
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").Value = ""
$ws.Range("H131").Value = 2575.5715
$ws.Range("I131").Value = 2575.5715
$ws.Range("K131").Value = 7726.7145
$ws.Range("M131").Value = -2686.7145
$ws.Range("H134").Value = 1423.3334
$ws.Range("I134").Value = 1423.3334
$ws.Range("K134").Value = 4270.0002
$ws.Range("M134").Value = 799.9997999999996
$ws.Range("H137").Value = 16668266
$ws.Range("I137").Value = 50001000
$ws.Range("J137").Value = 1899.5
$ws.Range("K137").Value = 150003000
$ws.Range("L137").Value = 5698.5
$ws.Range("M137").Value = -149997900
$ws.Range("N137").Value = -15898.5
$ws.Range("H138").Value = 1380
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""
# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 665.1111
$ws.Range("I2").Value = 922.4167
$ws.Range("K2").Value = 922.4167
$ws.Range("M2").Value = -809.4167
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = ""
$ws.Range("H38").Value = 24512
$ws.Range("I38").Value = 29000
$ws.Range("J38").Value = 20024
$ws.Range("K38").Value = 29000
$ws.Range("L38").Value = 20024
$ws.Range("M38").Value = -28537
$ws.Range("N38").Value = -20950
$ws.Range("H46").Value = 31000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 31000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 31000
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -31312
$ws.Range("H63").Value = 26050
$ws.Range("I63").Value = 26050
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 26050
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -25364
$ws.Range("N63").Value = ""
$ws.Range("H66").Value = 26050
$ws.Range("I66").Value = 26050
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 78150
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -74718
$ws.Range("N66").Value = ""
$ws.Range("H70").Value = 3193.8
$ws.Range("I70").Value = 3326.3333
$ws.Range("K70").Value = 3326.3333
$ws.Range("M70").Value = -3056.3333
$ws.Range("H73").Value = 3193.8
$ws.Range("I73").Value = 3326.3333
$ws.Range("K73").Value = 3326.3333
$ws.Range("M73").Value = -2390.3333
$ws.Range("H80").Value = 2936.7778
$ws.Range("I80").Value = 2490.2856
$ws.Range("J80").Value = 4499.5
$ws.Range("K80").Value = 2490.2856
$ws.Range("L80").Value = 4499.5
$ws.Range("M80").Value = -1492.2856
$ws.Range("N80").Value = -6495.5
$ws.Range("H83").Value = 2936.7778
$ws.Range("I83").Value = 2490.2856
$ws.Range("J83").Value = 4499.5
$ws.Range("K83").Value = 12451.428
$ws.Range("L83").Value = 22497.5
$ws.Range("M83").Value = -7459.428
$ws.Range("N83").Value = -32481.5
$ws.Range("H113").Value = 34261.477
$ws.Range("I113").Value = 38221.723
$ws.Range("K113").Value = 38221.723
$ws.Range("M113").Value = -36051.723
$ws.Range("H132").Value = 2317381.5
$ws.Range("J132").Value = 4942.3335
$ws.Range("L132").Value = 14827.0005
$ws.Range("N132").Value = -19887.0005
$ws.Range("H136").Value = 70747.125
$ws.Range("J136").Value = 70747.125
$ws.Range("L136").Value = 212241.375
$ws.Range("N136").Value = -217341.375
# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H20").Value = 29649.834
$ws.Range("I20").Value = 35966.332
$ws.Range("J20").Value = 23333.334
$ws.Range("K20").Value = 35966.332
$ws.Range("L20").Value = 23333.334
$ws.Range("M20").Value = -35740.332
$ws.Range("N20").Value = -23785.334
$ws.Range("H22").Value = 2512.4614
$ws.Range("J22").Value = 4000
$ws.Range("L22").Value = 4000
$ws.Range("N22").Value = -4590
$ws.Range("H27").Value = 2512.4614
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 4000
$ws.Range("N27").Value = -4214
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("H55").Value = 512.6
$ws.Range("I55").Value = 512.6
$ws.Range("K55").Value = 512.6
$ws.Range("M55").Value = -339.6
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = ""
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = ""
$ws.Range("H68").Value = 2594.5
$ws.Range("I68").Value = 3334.25
$ws.Range("K68").Value = 3334.25
$ws.Range("M68").Value = -2585.25
$ws.Range("H71").Value = 2594.5
$ws.Range("I71").Value = 3334.25
$ws.Range("K71").Value = 16671.25
$ws.Range("M71").Value = -12927.25
$ws.Range("H82").Value = 2567.3333
$ws.Range("I82").Value = 1351
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 1351
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -990
$ws.Range("N82").Value = -5722
$ws.Range("H85").Value = 2567.3333
$ws.Range("I85").Value = 1351
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 1351
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -103
$ws.Range("N85").Value = -7496
$ws.Range("H132").Value = 5003953.5
$ws.Range("I132").Value = 5105973
$ws.Range("K132").Value = 15317919
$ws.Range("M132").Value = -15315389
$ws.Range("H136").Value = 1960
$ws.Range("I136").Value = 1960
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5880
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3330
$ws.Range("N136").Value = ""
# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H69").Value = 22702.857
$ws.Range("J69").Value = 22702.857
$ws.Range("L69").Value = 22702.857
$ws.Range("N69").Value = -24200.857
$ws.Range("H72").Value = 22702.857
$ws.Range("J72").Value = 22702.857
$ws.Range("L72").Value = 68108.571
$ws.Range("N72").Value = -75596.571
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676
$ws.Range("H122").Value = 912.93335
$ws.Range("I122").Value = 912.93335
$ws.Range("K122").Value = 2738.80005
$ws.Range("M122").Value = -288.8000499999998
$ws.Range("H126").Value = 1588.6666
$ws.Range("I126").Value = 1159.8
$ws.Range("K126").Value = 3479.4
$ws.Range("M126").Value = -1009.4
$ws.Range("H132").Value = 10872174
$ws.Range("I132").Value = 13891381
$ws.Range("K132").Value = 41674143
$ws.Range("M132").Value = -41671613
